$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 58
$ws.Range("I2").Value = 147
$ws.Range("J2").Value = 647
$ws.Range("L2").Value = 177
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 114
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 67
$ws.Range("T2").Value = 123
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 1007
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 983
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 7
